# Added generic constraints for processes.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. scenarios sheet: tweak probabilities and update selection
# ---------------------------------------------------------------------------
$scenarios = $wb.Worksheets.Item("scenarios")
$scenarios.Activate()
$scenarios.Range("B2").Value = 0.3
$scenarios.Range("B4").Value = 0.2
$scenarios.Range("B2:B4").Select()

# ---------------------------------------------------------------------------
# 2. inflow sheet: zero out the first scenario's inflow values, update selection
# ---------------------------------------------------------------------------
$inflow = $wb.Worksheets.Item("inflow")
$inflow.Activate()
$inflow.Range("B2").Value = 0
$inflow.Range("C2").Value = 0
$inflow.Range("D2").Value = 0
$inflow.Range("B2").Select()

# ---------------------------------------------------------------------------
# 3. process_topology sheet: update selection only
# ---------------------------------------------------------------------------
$topology = $wb.Worksheets.Item("process_topology")
$topology.Activate()
$topology.Range("E2").Select()

# ---------------------------------------------------------------------------
# 4. processes sheet: update selection only (loses tabSelected later)
# ---------------------------------------------------------------------------
$processes = $wb.Worksheets.Item("processes")
$processes.Activate()
$processes.Range("H3:I3").Select()

# ---------------------------------------------------------------------------
# 5. New sheet: constraints
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$constraints = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$constraints.Name = "constraints"
$constraints.Activate()

$constraints.Range("A1").Value = "name"
$constraints.Range("B1").Value = "type"
$constraints.Range("A2").Value = "c1"
$constraints.Range("B2").Value = "eq"
$constraints.Range("I17").Select()

# ---------------------------------------------------------------------------
# 6. New sheet: gen_constraint
# ---------------------------------------------------------------------------
$genConstraint = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $constraints)
$genConstraint.Name = "gen_constraint"
$genConstraint.Activate()

$genConstraint.Range("A1").Value = "t"
$genConstraint.Range("B1").Value = "c1,ngchp,elc,s1"
$genConstraint.Range("C1").Value = "c1,ngchp,elc,s2"
$genConstraint.Range("D1").Value = "c1,ngchp,elc,s3"
$genConstraint.Range("E1").Value = "c1,ngchp,dh,s1"
$genConstraint.Range("F1").Value = "c1,ngchp,dh,s2"
$genConstraint.Range("G1").Value = "c1,ngchp,dh,s3"
$genConstraint.Range("H1").Value = "c1,s1"
$genConstraint.Range("I1").Value = "c1,s2"
$genConstraint.Range("J1").Value = "c1,s3"

# Column A: reuse the exact time-of-day fractions used by the other time
# series sheets (copy raw values to keep bit-identical doubles).
$fixedTs = $wb.Worksheets.Item("fixed_ts")
$timeVals = $fixedTs.Range("A2:A25").Value2
$genConstraint.Range("A2:A25").Value = $timeVals
$genConstraint.Range("A2:A25").NumberFormat = "h:mm"

$rowBC = New-Object 'object[,]' 24,9
for ($i = 0; $i -lt 24; $i++) {
    $rowBC[$i,0] = 1
    $rowBC[$i,1] = 1
    $rowBC[$i,2] = 1
    $rowBC[$i,3] = -0.8
    $rowBC[$i,4] = -0.8
    $rowBC[$i,5] = -0.8
    $rowBC[$i,6] = 0
    $rowBC[$i,7] = 0
    $rowBC[$i,8] = 0
}
$genConstraint.Range("B2:J25").Value = $rowBC

$genConstraint.Range("E4").Select()

# ---------------------------------------------------------------------------
# 7. nodes sheet: final active sheet/selection
# ---------------------------------------------------------------------------
$nodes = $wb.Worksheets.Item("nodes")
$nodes.Activate()
$nodes.Range("B10").Select()
